$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-03 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-04 Monday", 2) | Out-Null
$d.Content.Find.Execute("10÷3=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "78÷3=26, 0", 2) | Out-Null
$d.Content.Find.Execute("63÷8=7, 7", $true, $false, $false, $false, $false, $true, 1, $false, "48÷6=8, 0", 2) | Out-Null
$d.Content.Find.Execute("34÷2=17, 0", $true, $false, $false, $false, $false, $true, 1, $false, "23÷2=11, 1", 2) | Out-Null
$d.Content.Find.Execute("46÷3=15, 1", $true, $false, $false, $false, $false, $true, 1, $false, "24÷9=2, 6", 2) | Out-Null
$d.Content.Find.Execute("78÷9=8, 6", $true, $false, $false, $false, $false, $true, 1, $false, "37÷7=5, 2", 2) | Out-Null
$d.Content.Find.Execute("56÷5=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "57÷6=9, 3", 2) | Out-Null
$d.Content.Find.Execute("56÷7=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "78÷6=13, 0", 2) | Out-Null
$d.Content.Find.Execute("32÷6=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "35÷7=5, 0", 2) | Out-Null
$d.Content.Find.Execute("19÷4=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "64÷7=9, 1", 2) | Out-Null
$d.Content.Find.Execute("25÷6=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "72÷8=9, 0", 2) | Out-Null
$d.Content.Find.Execute("46÷7=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "62÷3=20, 2", 2) | Out-Null
$d.Content.Find.Execute("33÷6=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "45÷7=6, 3", 2) | Out-Null
$d.Content.Find.Execute("73÷4=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "17÷9=1, 8", 2) | Out-Null
$d.Content.Find.Execute("86÷3=28, 2", $true, $false, $false, $false, $false, $true, 1, $false, "34÷2=17, 0", 2) | Out-Null
$d.Content.Find.Execute("94÷2=47, 0", $true, $false, $false, $false, $false, $true, 1, $false, "24÷7=3, 3", 2) | Out-Null
$d.Content.Find.Execute("67÷3=22, 1", $true, $false, $false, $false, $false, $true, 1, $false, "97÷3=32, 1", 2) | Out-Null
$d.Content.Find.Execute("61÷8=7, 5", $true, $false, $false, $false, $false, $true, 1, $false, "52÷4=13, 0", 2) | Out-Null
$d.Content.Find.Execute("43÷2=21, 1", $true, $false, $false, $false, $false, $true, 1, $false, "86÷7=12, 2", 2) | Out-Null
$d.Content.Find.Execute("42÷6=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "95÷9=10, 5", 2) | Out-Null
$d.Content.Find.Execute("49÷8=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "67÷7=9, 4", 2) | Out-Null
$d.Content.Find.Execute("17÷6=2, 5", $true, $false, $false, $false, $false, $true, 1, $false, "17÷3=5, 2", 2) | Out-Null
$d.Content.Find.Execute("21÷6=3, 3", $true, $false, $false, $false, $false, $true, 1, $false, "38÷4=9, 2", 2) | Out-Null
$d.Content.Find.Execute("41÷6=6, 5", $true, $false, $false, $false, $false, $true, 1, $false, "21÷2=10, 1", 2) | Out-Null
$d.Content.Find.Execute("36÷6=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "73÷5=14, 3", 2) | Out-Null
$d.Content.Find.Execute("80÷9=8, 8", $true, $false, $false, $false, $false, $true, 1, $false, "88÷9=9, 7", 2) | Out-Null

Write-Output "Replacements applied"